$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = questions, Column B = answers.
# Row 1 keeps its taller height (28.8); rows 2, 6, 8 revert to default height.

$ws.Range("A1").Value = "Where can I find access review calendar?`nWhen the next review kick off?"
$ws.Range("B1").Value = "Here is the current access review calendar"

$ws.Range("A2").Value = "Where can I find application contact."
$ws.Range("B2").Value = "Here is the application contact information."

# NOTE: row 4's new question text is inserted into the shared-string table
# before row 3's (matches the target shared-string index order), so set A4
# ahead of A3.
$ws.Range("A4").Value = "What's entitlement risk level?"
$ws.Range("A3").Value = "How to reassign reviews?"

$ws.Range("B3").Value = "Here is the information for review reassignment."
$ws.Range("B4").Value = "Here is the information for High Risk Level."

$ws.Range("A5").Value = "How to update entitlement description?"
$ws.Range("B5").Value = "Here is the information regarding entitlement description."

$ws.Range("A6").Value = "How to update AD group to my application?"
$ws.Range("B6").Value = "Here is the information for AD groups."

$ws.Range("A7").Value = "How to get data query document?"
$ws.Range("B7").Value = "Here is the information for query scripts."

$ws.Range("A8").Value = "How to upload ACL file?"
$ws.Range("B8").Value = "Here is the information for ACL sFTP."

# Row heights: only row 1 keeps a custom (taller) height; rows 2, 6, 8 go back to default (auto) height.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(8).AutoFit()

# Column A widened to fit the new (longer) questions (target stored width ~64.109375
# character-units; the host snaps ColumnWidth to whole-pixel steps, so feed it the
# input that lands closest on that grid).
$ws.Columns.Item(1).ColumnWidth = 63.3

# Selection moves to A9 (just past the last used row).
$ws.Range("A9").Select()
